$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G21/H21 values
$ws.Range("G21").Value = 0.022
$ws.Range("H21").Value = 4

# Delete entire row 34 (Parapenaeus longirostris/PAPELON duplicate), shifting rows 35-38 up
$ws.Rows(34).Delete()
